$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "ont_demo" worksheet after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "ont_demo"

$ws.Range("A1").Value = "Especificacion"
$ws.Range("B1").Value = "Concepto/RegExp/Pair"
$ws.Range("C1").Value = "ClienteTexto_Campos"
$ws.Range("D1").Value = "ClienteTexto_Especificacion"

$ws.Range("A2").Value = "Discriminación > Sexual | Insulto"
$ws.Range("B2").Value = "puto"
$ws.Range("C2").Value = "texto"
$ws.Range("D2").Value = "TST_RechAuto_Insulto_SE_Normal"

$ws.Range("A3").Value = "Insulto"
$ws.Range("B3").Value = "boludo"
$ws.Range("C3").Value = "texto"
$ws.Range("D3").Value = "TST_ModMan_Insulto_SU_Normal"

$ws.Range("A1:D1").Font.Name = "Calibri"
$ws.Range("A1:D1").Font.Size = 10
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").Font.Color = 0

$ws.Range("A2:D3").Font.Name = "Calibri"
$ws.Range("A2:D3").Font.Size = 10
$ws.Range("A2:D3").Font.Bold = $false
$ws.Range("A2:D3").Font.Color = 0

# ---------------------------------------------------------------------------
# 2. test_otros: B2 becomes a real =FALSE() formula.
# ---------------------------------------------------------------------------
$wsOtros = $wb.Worksheets.Item("test_otros")
$wsOtros.Range("B2").Formula = "=FALSE()"

# ---------------------------------------------------------------------------
# 3. test_param: new _FilterDatabase name.
# ---------------------------------------------------------------------------
$wb.Names.Add("_xlnm._FilterDatabase_1", "=test_param!`$A`$1:`$L`$1")
